# Weekly update: a new price-report row is inserted at row 28 (pushing the
# existing rows 28-37 down to 29-38); the new row carries the latest
# "Poroto verde" quote for Terminal Hortofrutícola Agro Chillán.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 28; Excel shifts rows 28:37 down to 29:38 and
# keeps the existing formatting (e.g. the date style on column D).
$ws.Rows("28").Insert()

# Populate the newly inserted row 28 with this week's data.
$ws.Range("A28").Value = 7
$ws.Range("B28").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C28").Value = "Ñuble"
$ws.Range("D28").Value = 44524
$ws.Range("E28").Value = 16
$ws.Range("F28").Value = 100112031
$ws.Range("G28").Value = "Poroto verde"
$ws.Range("H28").Value = "Magnum"
$ws.Range("I28").Value = "Primera"
$ws.Range("J28").Value = 60
$ws.Range("K28").Value = 29000
$ws.Range("L28").Value = 30000
$ws.Range("M28").Value = 29500
$ws.Range("N28").Value = "$/saco 25 kilos"
$ws.Range("O28").Value = "Región Metropolitana"
$ws.Range("P28").Value = 1180
$ws.Range("Q28").Value = 25
$ws.Range("R28").Value = "Hortaliza"
